# Insert a new data row at row 483 on the active sheet.
# This shifts the existing rows 483:594 down to 484:595 and creates a
# brand-new row 483 with its own date/price data (same market/category
# metadata as the row that used to occupy row 483).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 483, pushing 483-594 -> 484-595.
$ws.Rows("483:483").Insert()

# Populate the newly inserted row 483 with the new record's data.
$ws.Range("A483").Value2 = 5
$ws.Range("B483").Value2 = "Macroferia Regional de Talca"
$ws.Range("C483").Value2 = "Maule"
$ws.Range("D483").Value2 = 45173
$ws.Range("E483").Value2 = 7
$ws.Range("F483").Value2 = 100114013
$ws.Range("G483").Value2 = "Zanahoria"
$ws.Range("H483").Value2 = "Sin especificar"
$ws.Range("I483").Value2 = "Primera"
$ws.Range("J483").Value2 = 600
$ws.Range("K483").Value2 = 5000
$ws.Range("L483").Value2 = 5000
$ws.Range("M483").Value2 = 5000
$ws.Range("N483").Value2 = "$/saco 20 kilos"
$ws.Range("O483").Value2 = "Región de Ñuble"
$ws.Range("P483").Value2 = 250
$ws.Range("Q483").Value2 = 20
$ws.Range("R483").Value2 = "Hortaliza"
